# [Kadastro App] Yeni kayit eklendi: 2977
# Adds the new work-record row (Kayit No 2977) to both the master
# "Kayitlar" log sheet and the per-district "Erdemli" sheet, which are
# kept in sync in this workbook.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "2977"
    B = "2025-09-10"
    C = "Erdemli"
    D = "1"
    E = "ÇAP"
    F = "CEMAL TİMUROĞLU (K.Teknisyeni)"
}

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find first empty row right below the existing data (row 44 here).
    $targetRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

    # Leading apostrophe keeps these text-looking numbers/dates stored as
    # text (matching every other row in the sheet) instead of letting
    # Excel auto-convert them to a number / date serial.
    $ws.Cells.Item($targetRow, 1).Value = "'" + $newRow.A
    $ws.Cells.Item($targetRow, 2).Value = "'" + $newRow.B
    $ws.Cells.Item($targetRow, 3).Value = $newRow.C
    $ws.Cells.Item($targetRow, 4).Value = "'" + $newRow.D
    $ws.Cells.Item($targetRow, 5).Value = $newRow.E
    $ws.Cells.Item($targetRow, 6).Value = $newRow.F
}
